$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Monthly Budget summary row
$ws.Range("F5").Value = 6311.458121275088
$ws.Range("G5").Value = 3800
$ws.Range("K5").Value = 44769.86551365722

# Other Debt amortization schedule (columns S:Y), rows 8-19 keep their
# row structure but get new computed values (loan paid off faster).
$ws.Range("U8").Value = 6311.458121275088
$ws.Range("V8").Value = 52177.49745530105
$ws.Range("W8").Value = 23560
$ws.Range("X8").Value = 75737.49745530105
$ws.Range("Y8").Value = 707822.502544699

$ws.Range("U9").Value = 6311.458121275088
$ws.Range("V9").Value = 53794.99987641539
$ws.Range("W9").Value = 21942.49757888567
$ws.Range("X9").Value = 75737.49745530105
$ws.Range("Y9").Value = 654027.5026682836

$ws.Range("U10").Value = 6311.458121275088
$ws.Range("V10").Value = 55462.64487258426
$ws.Range("W10").Value = 20274.85258271679
$ws.Range("X10").Value = 75737.49745530105
$ws.Range("Y10").Value = 598564.8577956994

$ws.Range("U11").Value = 6311.458121275088
$ws.Range("V11").Value = 57181.98686363437
$ws.Range("W11").Value = 18555.51059166668
$ws.Range("X11").Value = 75737.49745530105
$ws.Range("Y11").Value = 541382.870932065

$ws.Range("U12").Value = 6311.458121275088
$ws.Range("V12").Value = 58954.62845640704
$ws.Range("W12").Value = 16782.86899889401
$ws.Range("X12").Value = 75737.49745530105
$ws.Range("Y12").Value = 482428.2424756579

$ws.Range("U13").Value = 6311.458121275088
$ws.Range("V13").Value = 60782.22193855566
$ws.Range("W13").Value = 14955.27551674539
$ws.Range("X13").Value = 75737.49745530105
$ws.Range("Y13").Value = 421646.0205371022

$ws.Range("U14").Value = 6311.458121275088
$ws.Range("V14").Value = 62666.47081865089
$ws.Range("W14").Value = 13071.02663665017
$ws.Range("X14").Value = 75737.49745530105
$ws.Range("Y14").Value = 358979.5497184513

$ws.Range("U15").Value = 6311.458121275088
$ws.Range("V15").Value = 64609.13141402906
$ws.Range("W15").Value = 11128.36604127199
$ws.Range("X15").Value = 75737.49745530105
$ws.Range("Y15").Value = 294370.4183044223

$ws.Range("U16").Value = 6311.458121275088
$ws.Range("V16").Value = 66612.01448786397
$ws.Range("W16").Value = 9125.48296743709
$ws.Range("X16").Value = 75737.49745530105
$ws.Range("Y16").Value = 227758.4038165583

$ws.Range("U17").Value = 6311.458121275088
$ws.Range("V17").Value = 68676.98693698774
$ws.Range("W17").Value = 7060.510518313307
$ws.Range("X17").Value = 75737.49745530105
$ws.Range("Y17").Value = 159081.4168795706

$ws.Range("U18").Value = 6311.458121275088
$ws.Range("V18").Value = 70805.97353203436
$ws.Range("W18").Value = 4931.523923266688
$ws.Range("X18").Value = 75737.49745530105
$ws.Range("Y18").Value = 88275.44334753622

$ws.Range("U19").Value = 6311.458121275088
$ws.Range("V19").Value = 73000.95871152743
$ws.Range("W19").Value = 2736.538743773623
$ws.Range("X19").Value = 75737.49745530105
$ws.Range("Y19").Value = 15274.48463600878

# Row 20 becomes the "Total" row for the Other Debt schedule (loan is now
# paid off in 2035 instead of 2038), replacing what used to be row 24.
$ws.Range("T20").Clear()
$ws.Range("U20").Clear()
$ws.Range("Y20").Clear()
$ws.Range("S20").Value = "Total"
$ws.Range("V20").Value = 744725.5153639911
$ws.Range("W20").Value = 164124.4540996214
$ws.Range("X20").Value = 908849.9694636124

# Rows 21-23 no longer have Other Debt schedule data (loan paid off earlier).
$ws.Range("S21:Y21").Clear()
$ws.Range("S22:Y22").Clear()
$ws.Range("S23:Y23").Clear()

# Row 24's old "Total" line (S,V,W,X) is removed, superseded by row 20.
$ws.Range("S24").Clear()
$ws.Range("V24").Clear()
$ws.Range("W24").Clear()
$ws.Range("X24").Clear()
